$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Il23a"
$ws.Cells.Item(2,3).Value = "Il23r"
$ws.Cells.Item(2,4).Value = "MuSCs"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 1.270105666666667
$ws.Cells.Item(2,8).Value = 3.810317
$ws.Cells.Item(2,9).Value = 0.220415243968963
$ws.Cells.Item(2,10).Value = 0.220415243968963
$ws.Cells.Item(2,11).Value = 1
$ws.Cells.Item(2,12).Value = 0.3333333333333333
$ws.Cells.Item(2,13).Value = 0.04298466666666667
$ws.Cells.Item(2,14).Value = 0.128954
$ws.Cells.Item(2,15).Value = 1
$ws.Cells.Item(2,16).Value = 1
$ws.Cells.Item(2,17).Value = 0.05459506871311112
$ws.Cells.Item(2,18).Value = 0.491355618418
$ws.Cells.Item(2,19).Value = 0.220415243968963
$ws.Cells.Item(2,20).Value = 0.220415243968963

# Row 3
$ws.Cells.Item(3,1).Value = "FAPs"
$ws.Cells.Item(3,2).Value = "Il23a"
$ws.Cells.Item(3,3).Value = "Il23r"
$ws.Cells.Item(3,4).Value = "MuSCs"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 1.653878333333333
$ws.Cells.Item(3,8).Value = 4.961635
$ws.Cells.Item(3,9).Value = 0.287015486903044
$ws.Cells.Item(3,10).Value = 0.287015486903044
$ws.Cells.Item(3,11).Value = 1
$ws.Cells.Item(3,12).Value = 0.3333333333333333
$ws.Cells.Item(3,13).Value = 0.04298466666666667
$ws.Cells.Item(3,14).Value = 0.128954
$ws.Cells.Item(3,15).Value = 1
$ws.Cells.Item(3,16).Value = 1
$ws.Cells.Item(3,17).Value = 0.07109140886555557
$ws.Cells.Item(3,18).Value = 0.6398226797900001
$ws.Cells.Item(3,19).Value = 0.287015486903044
$ws.Cells.Item(3,20).Value = 0.287015486903044

# Row 4
$ws.Cells.Item(4,1).Value = "Inflammatory-Mac"
$ws.Cells.Item(4,2).Value = "Il23a"
$ws.Cells.Item(4,3).Value = "Il23r"
$ws.Cells.Item(4,4).Value = "MuSCs"
$ws.Cells.Item(4,5).Value = 2
$ws.Cells.Item(4,6).Value = 0.6666666666666666
$ws.Cells.Item(4,7).Value = 1.447887666666667
$ws.Cells.Item(4,8).Value = 4.343663
$ws.Cells.Item(4,9).Value = 0.2512676871409801
$ws.Cells.Item(4,10).Value = 0.2512676871409801
$ws.Cells.Item(4,11).Value = 1
$ws.Cells.Item(4,12).Value = 0.3333333333333333
$ws.Cells.Item(4,13).Value = 0.04298466666666667
$ws.Cells.Item(4,14).Value = 0.128954
$ws.Cells.Item(4,15).Value = 1
$ws.Cells.Item(4,16).Value = 1
$ws.Cells.Item(4,17).Value = 0.06223696872244445
$ws.Cells.Item(4,18).Value = 0.5601327185020001
$ws.Cells.Item(4,19).Value = 0.2512676871409801
$ws.Cells.Item(4,20).Value = 0.2512676871409801

# Row 5
$ws.Cells.Item(5,1).Value = "MuSCs"
$ws.Cells.Item(5,2).Value = "Il23a"
$ws.Cells.Item(5,3).Value = "Il23r"
$ws.Cells.Item(5,4).Value = "MuSCs"
$ws.Cells.Item(5,5).Value = 2
$ws.Cells.Item(5,6).Value = 0.6666666666666666
$ws.Cells.Item(5,7).Value = 0.324919
$ws.Cells.Item(5,8).Value = 0.974757
$ws.Cells.Item(5,9).Value = 0.05638672634467276
$ws.Cells.Item(5,10).Value = 0.05638672634467275
$ws.Cells.Item(5,11).Value = 1
$ws.Cells.Item(5,12).Value = 0.3333333333333333
$ws.Cells.Item(5,13).Value = 0.04298466666666667
$ws.Cells.Item(5,14).Value = 0.128954
$ws.Cells.Item(5,15).Value = 1
$ws.Cells.Item(5,16).Value = 1
$ws.Cells.Item(5,17).Value = 0.01396653490866667
$ws.Cells.Item(5,18).Value = 0.125698814178
$ws.Cells.Item(5,19).Value = 0.05638672634467276
$ws.Cells.Item(5,20).Value = 0.05638672634467275

# New row 6
$ws.Cells.Item(6,1).Value = "Resolving-Mac"
$ws.Cells.Item(6,2).Value = "Il23a"
$ws.Cells.Item(6,3).Value = "Il23r"
$ws.Cells.Item(6,4).Value = "MuSCs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 1.065540666666667
$ws.Cells.Item(6,8).Value = 3.196622
$ws.Cells.Item(6,9).Value = 0.1849148556423401
$ws.Cells.Item(6,10).Value = 0.1849148556423401
$ws.Cells.Item(6,11).Value = 1
$ws.Cells.Item(6,12).Value = 0.3333333333333333
$ws.Cells.Item(6,13).Value = 0.04298466666666667
$ws.Cells.Item(6,14).Value = 0.128954
$ws.Cells.Item(6,15).Value = 1
$ws.Cells.Item(6,16).Value = 1
$ws.Cells.Item(6,17).Value = 0.04580191037644445
$ws.Cells.Item(6,18).Value = 0.4122171933880001
$ws.Cells.Item(6,19).Value = 0.1849148556423401
$ws.Cells.Item(6,20).Value = 0.1849148556423401
